$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.808.88"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.091.88"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.86"
$ws.Range("E5").Value = "  -2.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.98"
$ws.Range("E6").Value = "  -1.39%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.083.34"
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("E10").Value = "  -3.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.34"
$ws.Range("E11").Value = "  -3.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.460"
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000228"
$ws.Range("E13").Value = "  +4.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.96"
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.590.09"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.761.28"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.090.71"
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "489.49"
$ws.Range("E20").Value = "  -3.60%  "
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.704"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.21"
$ws.Range("E23").Value = "  -0.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.01"
$ws.Range("E24").Value = "  +2.30%  "
$ws.Range("E25").Value = "  -1.21%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  -1.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.30"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.29"
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("E31").Value = "  -2.39%  "
$ws.Range("E32").Value = "  -0.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.41"
$ws.Range("E33").Value = "  -5.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "57.06"
$ws.Range("E34").Value = "  -2.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.49"
$ws.Range("E35").Value = "  +4.75%  "
$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "496.51"
$ws.Range("E36").Value = "  -6.45%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.08"
$ws.Range("E37").Value = "  +1.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.319.48"
$ws.Range("E38").Value = "  +7.85%  "
$ws.Range("E39").Value = "  -3.56%  "
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("E41").Value = "  -2.40%  "
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("E43").Value = "  -2.28%  "
$ws.Range("E44").Value = "  +2.37%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.09"
$ws.Range("E46").Value = "  +1.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0543"
$ws.Range("E47").Value = "  +6.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.17"
$ws.Range("E48").Value = "  +3.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.28"
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("E50").Value = "  +2.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.34"
$ws.Range("E51").Value = "  -14.36%  "
